$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.777.64'
$ws.Range("E2").Value = '  -0.51%  '
$ws.Range("D3").Value = '3.400.88'
$ws.Range("E3").Value = '  -0.73%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '407.11'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '127.96'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.627'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +5.88%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.723'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.54%  '
$ws.Range("E10").Value = '  +14.93%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.11'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.38%  '
$ws.Range("E12").Value = '  -0.69%  '
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '3.953.66'
$ws.Range("E13").Value = '  -0.22%  '
$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000212'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +60.75%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.85'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.74%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.75'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.16%  '
$ws.Range("D17").Value = '3.427.61'
$ws.Range("E17").Value = '  +0.51%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.00'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +8.62%  '
$ws.Range("E19").Value = '  +3.81%  '
$ws.Range("D20").Value = '61.709.53'
$ws.Range("E20").Value = '  -0.68%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '403.93'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +28.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '88.98'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.65%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.16'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.57%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.01'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.52%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.21'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.64%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '32.53'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +9.27%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.79'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.27%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.55'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.19%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.59'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.46%  '
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.118'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.39%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.68'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.31%  '
$ws.Range("E32").Value = '  -1.99%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.78'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.37%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '42.87'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.17%  '
$ws.Range("E35").Value = '  +0.69%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0493'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.31%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '53.86'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.23%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.34'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.86%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.132'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.60%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.90'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.94%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.308'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.66%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '139.81'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.50%  '
$ws.Range("E44").Value = '  -2.37%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.02'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.81%  '
$ws.Range("E46").Value = '  +8.42%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.52'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.83%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '21.59'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.05%  '
$ws.Range("D49").Value = '2.106.02'
$ws.Range("E49").Value = '  -1.18%  '
$ws.Range("E50").Value = '  +4.21%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.130'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +15.32%  '
